$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the placeholder empty cells C6:M6 (their content moves to the new row 7)
$ws.Range("C6:M6").ClearContents()

# Add new row 7 with data (duplicate of rows 4/5 pattern)
$ws.Range("A7").Value = 2488
$ws.Range("B7").Value = "Test Ringover (NO TOCAR)"
$ws.Range("C7").Value = "Estructura coplanar NOVOTEGRA"
$ws.Range("D7").Value = "'1"
$ws.Range("E7").Value = "'"
$ws.Range("F7").Value = "'"
$ws.Range("G7").Value = "'"
$ws.Range("H7").Value = "Inversor híbrido monofásico SUN-6k-SG05LP1-EU"
$ws.Range("I7").Value = "'1"
$ws.Range("J7").Value = "BATERÍA LITIO SIGEN ENERGY SIGENSTOR 10,0KW"
$ws.Range("K7").Value = "'3"
$ws.Range("L7").Value = "'"
$ws.Range("M7").Value = "Sí"
$ws.Range("N7").Value = "2024-01-03T10:49:29.104Z"
